$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 = duplicate of row 4 ("15 only ...")
$ws.Range("A5").Value = "15 only address"
$ws.Range("B5").Value = "15 only city"
$ws.Range("C5").Value = " "
$ws.Range("D5").Value = "15 only first"
$ws.Range("E5").Value = "15 only last"
$ws.Range("F5").Value = " "
$ws.Range("G5").Value = "{{address}}"
$ws.Range("H5").Value = "{{address}}"

# Row 6 = duplicate of row 3 ("15 ...")
$ws.Range("A6").Value = "15 address"
$ws.Range("B6").Value = "15 HCM city"
$ws.Range("C6").Value = " "
$ws.Range("D6").Value = "15 silicon first"
$ws.Range("E6").Value = "15 silicon last"
$ws.Range("F6").Value = "{{ip_address}}"
# Leading apostrophe forces these numeric-looking strings to stay text
# (so leading zeros / exact digits are preserved), then restore the
# default "Normal" style so no extra quote-prefix formatting lingers.
$ws.Range("G6").Value = "'1515"
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Value = "'151515"
$ws.Range("H6").Style = "Normal"

# Row 7 = duplicate of row 2 ("f03 ...")
$ws.Range("A7").Value = "f03 address"
$ws.Range("B7").Value = "03 city"
$ws.Range("C7").Value = " "
$ws.Range("D7").Value = "f03 first"
$ws.Range("E7").Value = "f03 last"
$ws.Range("F7").Value = " "
$ws.Range("G7").Value = "'0303"
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value = "'030303"
$ws.Range("H7").Style = "Normal"
